# The deck's slide-master theme ("Integral") is recoloured to use the
# stock "Office Theme" palette. dk1/lt1 (black/white) are identical in
# both themes, so only dk2, lt2 and the six accents plus the two
# hyperlink colors need to change; font scheme and format scheme are
# already identical between the two themes.
#
# PowerPoint's ThemeColorScheme.RGB setter (like VBA's RGB()) takes an
# OLE/BGR-packed integer, so each target hex color is converted with
# R | (G << 8) | (B << 16).

function BGR($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Index order (matches <a:clrScheme>): 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$tcs.Item(3).RGB  = BGR("44546A")  # dk2
$tcs.Item(4).RGB  = BGR("E7E6E6")  # lt2
$tcs.Item(5).RGB  = BGR("5B9BD5")  # accent1
$tcs.Item(6).RGB  = BGR("ED7D31")  # accent2
$tcs.Item(7).RGB  = BGR("A5A5A5")  # accent3
$tcs.Item(8).RGB  = BGR("FFC000")  # accent4
$tcs.Item(9).RGB  = BGR("4472C4")  # accent5
$tcs.Item(10).RGB = BGR("70AD47")  # accent6
$tcs.Item(11).RGB = BGR("0563C1")  # hlink
$tcs.Item(12).RGB = BGR("954F72")  # folHlink
